$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from row 2 to row 72: increment date serial
# from 45177 (2023-09-08) to 45178 (2023-09-09) for each data row.
for ($row = 2; $row -le 72; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
